# Incluye datos de SQL Server
# Adds a new "SQLServer" worksheet, populated with the same crop/area data as
# the other two sheets (sourced, in the original workbook, from a SQL Server
# connection via an .odc file), wraps it in a table, and defines the hidden
# workbook-scoped name that Excel creates for a query-table range.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet after the existing ones -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "SQLServer"

# --- 2. Populate it with the query results -----------------------------
$data = @(
  @("Cultivo", "Superficie"),
  @("Trigo", 53994657),
  @("Arroz en cáscara", 614453),
  @("Cebada", 29283399),
  @("Maíz", 13255122),
  @("Centeno", 8555737),
  @("Avena", 8096907),
  @("Mijo", 1029285),
  @("Sorgo", 154545),
  @("Alforfón", 1538326),
  @("Triticale", 1865252),
  @("Alpiste", 2207),
  @("Cereales mezclados", 1770317),
  @("Cereales nep", 114332),
  @("Patatas", 9144628)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# --- 3. Turn the range into the query table's table/list object --------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:B15"), $null, 1)
$lo.Name = "Tabla_WIN7DB_FAO_Cultivos"
$lo.TableStyle = "TableStyleMedium2"

# --- 4. Recreate the hidden, sheet-scoped defined name Excel keeps for --
#        the query table's range.
$name = $ws.Names.Add("WIN7DB_FAO_Cultivos", "=SQLServer!`$A`$1:`$B`$15")
$name.Visible = $false
